$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E (F unchanged); G is the row sum B+C+D+E
$data = @{
    2 = @(0.127881588408715, 0.3127903958511391, 0.1575252929769615, 0.496779210170732)
    3 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732)
    4 = @(0.01514828764759746, 0.002777888934908601, 0.1575252929769615, 0.496779210170732)
    5 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732)
    6 = @(0.127881588408715, 0.002777888934908601, 3.900430680208489, 8.660232485948974)
    7 = @(0.04763786555579896, 0.04240448674262143, 0.8054896365839992, 0.496779210170732)
    8 = @(1.459612070389937, 0.3127903958511391, 0.8054896365839992, 0.496779210170732)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $sum = $vals[0] + $vals[1] + $vals[2] + $vals[3]
    $ws.Cells.Item($row, 7).Value = $sum
}
